$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113 (this shifts rows 113-123 down to 114-124,
# and updates all existing merged cells / row heights / styles automatically).
$ws.Rows("113:113").Insert()

# Re-apply the bottom border used throughout the data rows so the new
# row reuses the same cell-style records as its neighbours instead of
# a border-less clone.
$rng = $ws.Range("A113:N113")
$rng.Borders.Item(9).LineStyle = 1
$rng.Borders.Item(9).Color = 13882323

# Populate the new product row: "ماسك جلسات اطفال"
$ws.Cells.Item(113, 1).Value = 110
$ws.Cells.Item(113, 2).Value = "ماسك جلسات اطفال"
$ws.Cells.Item(113, 8).Value = "-1:0"
$ws.Cells.Item(113, 12).Value = 20
$ws.Cells.Item(113, 14).Value = 1

# Recreate the merges for the newly-inserted row (Insert() does not
# carry merges onto a brand new blank row).
$ws.Range("B113:G113").Merge()
$ws.Range("H113:K113").Merge()
$ws.Range("L113:M113").Merge()

# Update the running total to include the new row's price (20).
$ws.Cells.Item(123, 11).Value = 6906.44

# Renumber the sequential index column (A) for every row pushed down by
# the insert, so the running count (1, 2, 3, ...) stays contiguous.
for ($r = 114; $r -le 122; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

